$p = $ppt.ActivePresentation

# --- Slide 1: quote text -> "Itghurls" ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = '"Itghurls"'

# --- Slide 2: tagline run (2nd run, after the soft line break) ---
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Characters(17, $tr2.Length - 16).Text = "PitchItup - a platform as a service for startup pitch deck."

# --- Slide 3: two bullet runs separated by a soft line break ---
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(3).TextFrame.TextRange
# replace the later run first so the earlier run's offsets stay valid
$tr3.Characters(21, 13).Text = "-  having hard time finding templates"
$tr3.Characters(1, 19).Text = "- Having hard time creating pitch deck"

# --- Slide 4: tagline run (2nd run, after the soft line break) ---
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Characters(13, $tr4.Length - 12).Text = "PitchItup - a platform as a service for startup pitch deck."

# --- Slide 5: "Benefits" bullet and "Features" bullet ---
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(3).TextFrame.TextRange.Text = "- Befmwkniwjk"
$s5.Shapes.Item(5).TextFrame.TextRange.Text = "- ijbwyuvew7beuh"
